# Updated symbol list on Wed Jan 18 10:57:56 UTC 2023 with GitHub Actions
#
# Refreshes Price/Volume(1h) figures for existing coin rows, and rotates the
# coin list for rows 16-35 so that "UpBots" (previously the last entry, row
# 35) now appears at row 16, with every other coin in that block shifting
# down by one row (their Price/Volume values follow along with a small
# refresh). Values that look numeric or percentage-like are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# workbook's existing text-based layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.90"
$ws.Range("E2").Value = "'0.07%"
$ws.Range("D3").Value = "'32.38"
$ws.Range("E3").Value = "'2.52%"
$ws.Range("D4").Value = "'4.946"
$ws.Range("E4").Value = "'-2.93%"
$ws.Range("D5").Value = "'0.07628"
$ws.Range("E5").Value = "'-2.29%"
$ws.Range("D6").Value = "'1.926"
$ws.Range("E6").Value = "'-14.12%"
$ws.Range("D7").Value = "'7.835"
$ws.Range("E7").Value = "'0.44%"
$ws.Range("D8").Value = "'0.9182"
$ws.Range("E8").Value = "'0.39%"
$ws.Range("D9").Value = "'0.1749"
$ws.Range("E9").Value = "'0.01%"
$ws.Range("D10").Value = "'0.07750"
$ws.Range("E10").Value = "'3.21%"
$ws.Range("D11").Value = "'0.08510"
$ws.Range("E11").Value = "'-5.40%"
$ws.Range("D12").Value = "'0.03195"
$ws.Range("E12").Value = "'3.16%"
$ws.Range("D13").Value = "'0.09999"
$ws.Range("E13").Value = "'-0.07%"
$ws.Range("D14").Value = "'0.001510"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("D15").Value = "'0.005930"
$ws.Range("E15").Value = "'-0.33%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007498"
$ws.Range("E16").Value = "'2,116.77%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.464"
$ws.Range("E17").Value = "'0.14%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.800"
$ws.Range("E18").Value = "'-0.83%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.152"
$ws.Range("E19").Value = "'-4.32%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3350"
$ws.Range("E20").Value = "'1.84%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1326"
$ws.Range("E21").Value = "'-0.85%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'4.273"
$ws.Range("E22").Value = "'5.71%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1992"
$ws.Range("E23").Value = "'9.58%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04520"
$ws.Range("E24").Value = "'-1.68%"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-2.24%"
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").Value = "'0.004388"
$ws.Range("E26").Value = "'-1.64%"
$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").Value = "'0.0001252"
$ws.Range("E27").Value = "'0.20%"
$ws.Range("B28").Value = "Spectre.aiUtilityToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("B29").Value = "LegolasExchange"
$ws.Range("C29").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("B30").Value = "BitZToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("B31").Value = "Birake"
$ws.Range("C31").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("B32").Value = "NashExchange"
$ws.Range("C32").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("B33").Value = "AAXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("B34").Value = "CenX"
$ws.Range("C34").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("B35").Value = "BNIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("D39").Value = "'0.01700"
$ws.Range("E39").Value = "'-4.34%"
$ws.Range("D40").Value = "'0.04688"
$ws.Range("E40").Value = "'-1.82%"
$ws.Range("D41").Value = "'0.007488"
$ws.Range("E41").Value = "'-0.91%"
$ws.Range("E42").Value = "'-0.64%"
$ws.Range("E43").Value = "'6.60%"
$ws.Range("D44").Value = "'0.01055"
$ws.Range("E44").Value = "'3.38%"
$ws.Range("D45").Value = "'0.00006259"
$ws.Range("E45").Value = "'0.85%"
$ws.Range("E46").Value = "'0.13%"
$ws.Range("D47").Value = "'0.8234"
$ws.Range("E47").Value = "'10.50%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.13%"
